# PFA_eksklusionsliste.xlsx maintenance edit.
#
# The list had accumulated a stray duplicate entry for the South-Korean
# defence company "S&T Holdings Co. Ltd." in row 92 (with inconsistent
# trailing whitespace compared to the other rows), sitting right above the
# correctly-formatted "Safran Group" entry further down the list. This
# edit removes that duplicate row, which shifts every following row up by
# one (rows 93-99 become 92-98) and shrinks the used range / shared
# string table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the stray duplicate row first (mirrors selecting a whole row in
# the Excel UI before deleting it).
$ws.Range("A92:XFD92").Select() | Out-Null

# Delete the row entirely - remaining rows shift up, and the sheet
# dimension / shared-string table shrink accordingly.
$ws.Rows.Item(92).Delete() | Out-Null
